# W1_V2 3DFormulas.xlsx - "week 1 of 2nd course"
#
# The "HR Q3" sheet is a roll-up of the three per-person expense sheets
# (Sean, Uma, Carlos). Its C/D/E detail cells (rows 7-13, 17-20, 24-27)
# were left blank; this fills them in with 3-D SUM formulas that total
# the corresponding cell across all three person sheets, e.g.
#   =SUM(Sean:Carlos!C7)
# which then flow up into the existing subtotal / grand-total formulas
# already present on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HR Q3")

$rows = @(7,8,9,10,11,12,13,17,18,19,20,24,25,26,27)
$cols = @("C","D","E")

foreach ($r in $rows) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $ws.Range($addr).Formula = "=SUM(Sean:Carlos!$addr)"
    }
}

# Reflect the author's final selection/scroll position on the HR Q3 tab.
$ws.Activate()
$ws.Range("C24:E27").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
